# Rename the two worksheets (ML/colorimetric-assay naming update) and
# switch the active tab from the first sheet to the second, matching the
# author's re-pointing of this workbook from the "machine_learning" data
# folder (rep3/rep4) to the "Dehalogenation_colorimetric_assays/Experiment1"
# folder (rep1/rep2).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "rep1_linearized_fluoride"
$ws2.Name = "rep2_linearized_fluoride"

# Make the second worksheet the active tab (moves tabSelected from sheet 1
# to sheet 2, and sets the workbook's activeTab to index 1).
$ws2.Activate()
